# Update the "F" column (views/likes count) numbers on several sheets,
# reflecting refreshed data pulled from the source at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1890
$ws1.Range("F3").Value = 1518
$ws1.Range("F5").Value = 775
$ws1.Range("F6").Value = 13355
$ws1.Range("F7").Value = 13221
$ws1.Range("F8").Value = 1021
$ws1.Range("F9").Value = 776
$ws1.Range("F11").Value = 564
$ws1.Range("F14").Value = 2098
$ws1.Range("F15").Value = 65
$ws1.Range("F17").Value = 75
$ws1.Range("F19").Value = 397
$ws1.Range("F20").Value = 265
$ws1.Range("F21").Value = 291
$ws1.Range("F23").Value = 762

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 131

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 49

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1890
$ws4.Range("F4").Value = 1518
$ws4.Range("F7").Value = 775
$ws4.Range("F8").Value = 13355
$ws4.Range("F9").Value = 13221
$ws4.Range("F10").Value = 1021
$ws4.Range("F11").Value = 776
$ws4.Range("F13").Value = 564
$ws4.Range("F18").Value = 2098
$ws4.Range("F19").Value = 65
$ws4.Range("F21").Value = 75
$ws4.Range("F22").Value = 131
$ws4.Range("F25").Value = 49
$ws4.Range("F26").Value = 397
$ws4.Range("F27").Value = 265
$ws4.Range("F28").Value = 291
$ws4.Range("F30").Value = 762
